$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 10.142424
$ws.Range("B4").Value = 9999

$ws.Range("B5").Select()
